$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Gremlin', ['Token Creature — Gremlin', '2/2', 'Energy Reserve', 'Card', '(Place your energy counters in this area.)'])"

$ws.Range("A3:A7").ClearContents()
